# Apply numeric tweaks to the NetMigration values and update the active
# cell selection, matching the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (NetMigration) values
$ws.Range("C3").Value  = 100000
$ws.Range("C4").Value  = 200000
$ws.Range("C5").Value  = 275000
$ws.Range("C6").Value  = 325000
$ws.Range("C7").Value  = 375000

$ws.Range("C9").Value  = 25000
$ws.Range("C10").Value = 30000
$ws.Range("C11").Value = 33000

$ws.Range("C13").Value = 36000

$ws.Range("C17").Value = 20000
$ws.Range("C18").Value = 25000
$ws.Range("C19").Value = 30000
$ws.Range("C20").Value = 10000
$ws.Range("C21").Value = 20000
$ws.Range("C22").Value = 25000
$ws.Range("C23").Value = 30000
$ws.Range("C24").Value = 35000
$ws.Range("C25").Value = 33000

# Move the active cell / selection to E18 (as recorded in the saved file)
$ws.Range("E18").Select()
